$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "307.16"
Set-TextValue $ws.Range("E2") "0.82%"
Set-TextValue $ws.Range("D3") "36.52"
Set-TextValue $ws.Range("E3") "3.25%"
Set-TextValue $ws.Range("D4") "5.067"
Set-TextValue $ws.Range("E4") "-0.54%"
Set-TextValue $ws.Range("D5") "0.08195"
Set-TextValue $ws.Range("E5") "2.43%"
Set-TextValue $ws.Range("D6") "1.966"
Set-TextValue $ws.Range("E6") "0.78%"
Set-TextValue $ws.Range("D7") "4.085"
Set-TextValue $ws.Range("E7") "-0.95%"
Set-TextValue $ws.Range("D8") "7.829"
Set-TextValue $ws.Range("E8") "-0.46%"
Set-TextValue $ws.Range("D9") "0.9347"
Set-TextValue $ws.Range("E9") "1.24%"
Set-TextValue $ws.Range("D10") "0.1477"
Set-TextValue $ws.Range("E10") "35.39%"
Set-TextValue $ws.Range("E11") "2.30%"
Set-TextValue $ws.Range("D12") "0.09128"
Set-TextValue $ws.Range("E12") "-3.19%"
Set-TextValue $ws.Range("D13") "0.03519"
Set-TextValue $ws.Range("E13") "-3.36%"
Set-TextValue $ws.Range("D14") "0.09812"
Set-TextValue $ws.Range("E14") "-0.96%"
Set-TextValue $ws.Range("D15") "0.001413"
Set-TextValue $ws.Range("E15") "0.11%"
Set-TextValue $ws.Range("D16") "0.005773"
Set-TextValue $ws.Range("E16") "-0.29%"
Set-TextValue $ws.Range("D17") "3.524"
Set-TextValue $ws.Range("E17") "2.10%"
Set-TextValue $ws.Range("E18") "5.27%"
Set-TextValue $ws.Range("D19") "0.3423"
Set-TextValue $ws.Range("E19") "0.07%"
Set-TextValue $ws.Range("D20") "0.1295"
Set-TextValue $ws.Range("E20") "-1.97%"
Set-TextValue $ws.Range("D21") "5.050"
Set-TextValue $ws.Range("E21") "-1.01%"
Set-TextValue $ws.Range("D22") "0.2394"
Set-TextValue $ws.Range("E22") "9.00%"
Set-TextValue $ws.Range("D23") "0.04498"
Set-TextValue $ws.Range("E23") "-0.69%"
Set-TextValue $ws.Range("D24") "0.001208"
Set-TextValue $ws.Range("E24") "-1.50%"
Set-TextValue $ws.Range("D25") "0.004906"
Set-TextValue $ws.Range("E25") "4.63%"
Set-TextValue $ws.Range("D26") "0.0001227"
Set-TextValue $ws.Range("E26") "-2.02%"
Set-TextValue $ws.Range("D27") "0.0004430"
Set-TextValue $ws.Range("E27") "-0.63%"
Set-TextValue $ws.Range("D39") "0.01988"
Set-TextValue $ws.Range("E39") "4.72%"
Set-TextValue $ws.Range("D40") "0.04859"
Set-TextValue $ws.Range("E40") "2.51%"
Set-TextValue $ws.Range("D41") "0.01105"
Set-TextValue $ws.Range("E41") "14.34%"
Set-TextValue $ws.Range("D42") "0.007541"
Set-TextValue $ws.Range("E42") "-0.04%"
Set-TextValue $ws.Range("D43") "0.1381"
Set-TextValue $ws.Range("E43") "3.32%"
Set-TextValue $ws.Range("D44") "0.002072"
Set-TextValue $ws.Range("E44") "-2.20%"
Set-TextValue $ws.Range("D45") "0.01090"
Set-TextValue $ws.Range("E45") "-3.49%"
Set-TextValue $ws.Range("D46") "0.00006109"
Set-TextValue $ws.Range("E46") "-4.16%"
Set-TextValue $ws.Range("E47") "-0.50%"
Set-TextValue $ws.Range("E48") "0.44%"
Set-TextValue $ws.Range("E49") "-8.87%"
Set-TextValue $ws.Range("D50") "0.00002095"
Set-TextValue $ws.Range("E50") "-0.50%"
Set-TextValue $ws.Range("D51") "0.0001995"
Set-TextValue $ws.Range("E51") "-0.50%"
